# Weekly data refresh for "Hortaliza, Vega Monumental Concepción - Perejil".
#
# A new week (2021-10-21, serial 44490) is added at the top of the data
# block (right after the two fixed "header" weeks kept in rows 2-9), which
# pushes every existing weekly pair of rows (Primera/Segunda) down by two
# rows. We reproduce that by duplicating rows 10:11 (the first data rows
# after the fixed block) and inserting the copy in their own place, which
# shifts the original rows 10:11 (and everything below) down by two rows.
# Afterwards we overwrite the date in the freshly inserted rows 10:11 with
# the new week's date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("10:11").Copy()
$ws.Rows("10:11").Insert()

$ws.Cells.Item(10, 4).Value = 44490
$ws.Cells.Item(11, 4).Value = 44490
